$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --------------------------------------------------------------------
# 1. Relocate existing rows downward to make room for the new
#    "Note on Sequence Co-ordinates" section (rows 7-8).
#    Processed from the bottom (highest target row) upward so that a
#    source row is never clobbered before it has been read.
# --------------------------------------------------------------------

$ws.Range("B21").Cut($ws.Range("B24"))
$ws.Range("B21").Clear()

$ws.Range("B20").Cut($ws.Range("B23"))
$ws.Range("B20").Clear()

$ws.Range("B18").Cut($ws.Range("B21"))
$ws.Range("B18").Clear()

$ws.Range("B16").Cut($ws.Range("B19"))
$ws.Range("B16").Clear()

$ws.Range("B13").Cut($ws.Range("B16"))
$ws.Range("B13").Clear()

$ws.Range("B11").Cut($ws.Range("B14"))
$ws.Range("B11").Clear()

$ws.Range("B8:C8").Cut($ws.Range("B11:C11"))
$ws.Range("B8:C8").Clear()

# --------------------------------------------------------------------
# 2. Add the new "Note on Sequence Co-ordinates" heading (row 7),
#    styled like the italic "Explanation of fields" heading that now
#    lives in B23.
# --------------------------------------------------------------------

$ws.Range("B23").Copy()
$ws.Range("B7").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B7").Value = "Note on Sequence Co-ordinates"

# --------------------------------------------------------------------
# 3. Add the explanatory note text in a merged, wrapped cell B8:L8.
# --------------------------------------------------------------------

$noteCell = $ws.Range("B8")
$noteCell.WrapText = $true
$noteCell.HorizontalAlignment = -4131 # xlLeft
$noteCell.Copy()
$ws.Range("C8:L8").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B8:L8").Merge()
$noteCell.Value = "Nucleotides within the inferred sequence are numbered from 1 - that is, a sequence of length n `nshould be numbered from 1 to n. These numbers are used as co-ordinates to identify the start `nand end nucleotide of defined regions."
$ws.Rows.Item(8).RowHeight = 48.75

# --------------------------------------------------------------------
# 4. Tidy up sheet-level bookkeeping to match the updated layout.
# --------------------------------------------------------------------

$ws.Range("A3:L24").Select()
$ws.Range("C14").Select()
